# Applies the sprint-log entries for "Week 24" to the Invulblad worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Numeric week number entered in E1
$ws.Range("E1").Value = 24

# Text entries entered in column C for the various team members / days
$ws.Range("C2").Value  = "Asylum"
$ws.Range("C4").Value  = "``Michiel"
$ws.Range("C6").Value  = "Sprint 7 Models Fixen , en Cutlery"
$ws.Range("C7").Value  = "Puzzle Script "
$ws.Range("C9").Value  = "Planning 7 , Longchairs ,smallpaper , CSI , Mirror"
$ws.Range("C10").Value = "Textures env , Lvl design . Implementing"
$ws.Range("C12").Value = "Kitchen Models"
$ws.Range("C14").Value = "Laundry room"
$ws.Range("C15").Value = "Wat ik op planning zet , en robin helpen"
$ws.Range("C17").Value = "Luxechair , Plants , Old tv , Paperstacks , Newspaper"
$ws.Range("C18").Value = "Wat ik op de planning zet , Bezig met lvl"
$ws.Range("C20").Value = "Wat ik Op de planning zet "
$ws.Range("C22").Value = "robbin echt heeeeeeeeeeel beidehand"
$ws.Range("C23").Value = "NOOOUUUU nee. Naast het feit dat tie geen internet heeft"
$ws.Range("C24").Value = "De chick is dood"
$ws.Range("C25").Value = "Ziek"
$ws.Range("C28").Value = "Helemaal mooi"
$ws.Range("C29").Value = "Planning word gevolged"
$ws.Range("C31").Value = "Assets sneller af"
